$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = ""
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H138").Value = 4295.1724
$ws.Range("I138").Value = 3228
$ws.Range("J138").Value = 6666.6665
$ws.Range("K138").Value = 9684
$ws.Range("L138").Value = 19999.9995
$ws.Range("M138").Value = -4544
$ws.Range("N138").Value = -30279.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3679043
$ws.Range("I2").Value = 2940
$ws.Range("J2").Value = 29411764
$ws.Range("K2").Value = 2940
$ws.Range("L2").Value = 29411764
$ws.Range("M2").Value = -2827
$ws.Range("N2").Value = -29411990
$ws.Range("H3").Value = 41669196
$ws.Range("H5").Value = 147
$ws.Range("I5").Value = 138.16667
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 138.16667
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -26.16667000000001
$ws.Range("N5").Value = -424
$ws.Range("H32").Value = 7967.5835
$ws.Range("I32").Value = 4499.783
$ws.Range("J32").Value = 25306.584
$ws.Range("K32").Value = 4499.783
$ws.Range("L32").Value = 25306.584
$ws.Range("M32").Value = -4212.783
$ws.Range("N32").Value = -25880.584
$ws.Range("H45").Value = 64173.062
$ws.Range("I45").Value = 92377.73
$ws.Range("J45").Value = 2122.8
$ws.Range("K45").Value = 92377.73
$ws.Range("L45").Value = 2122.8
$ws.Range("M45").Value = -92000.73
$ws.Range("N45").Value = -2876.8
$ws.Range("H61").Value = 6291379
$ws.Range("I61").Value = 9525342
$ws.Range("J61").Value = 3118.5
$ws.Range("K61").Value = 9525342
$ws.Range("L61").Value = 3118.5
$ws.Range("M61").Value = -9525130
$ws.Range("N61").Value = -3542.5
$ws.Range("H102").Value = 2044
$ws.Range("I102").Value = 2082.5
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 2082.5
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -460.5
$ws.Range("N102").Value = -5244
$ws.Range("H116").Value = 3679043
$ws.Range("I116").Value = 2940
$ws.Range("J116").Value = 29411764
$ws.Range("K116").Value = 2940
$ws.Range("L116").Value = 29411764
$ws.Range("M116").Value = -646
$ws.Range("N116").Value = -29416352
$ws.Range("H122").Value = 1928.909
$ws.Range("I122").Value = 1610.9546
$ws.Range("J122").Value = 2564.818
$ws.Range("K122").Value = 4832.8638
$ws.Range("L122").Value = 7694.454000000001
$ws.Range("M122").Value = -2382.8638
$ws.Range("N122").Value = -12594.454
$ws.Range("H136").Value = 6291379
$ws.Range("I136").Value = 9525342
$ws.Range("J136").Value = 3118.5
$ws.Range("K136").Value = 28576026
$ws.Range("L136").Value = 9355.5
$ws.Range("M136").Value = -28573476
$ws.Range("N136").Value = -14455.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3679043
$ws.Range("I3").Value = 2940
$ws.Range("J3").Value = 29411764
$ws.Range("K3").Value = 2940
$ws.Range("L3").Value = 29411764
$ws.Range("M3").Value = -2826
$ws.Range("N3").Value = -29411992
$ws.Range("H4").Value = 147
$ws.Range("I4").Value = 138.16667
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 138.16667
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -23.16667000000001
$ws.Range("N4").Value = -430
$ws.Range("H20").Value = 2666.8
$ws.Range("I20").Value = 2814.25
$ws.Range("J20").Value = 2498.2856
$ws.Range("K20").Value = 2814.25
$ws.Range("L20").Value = 2498.2856
$ws.Range("M20").Value = -2567.25
$ws.Range("N20").Value = -2992.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 21500
$ws.Range("I6").Value = 40000
$ws.Range("K6").Value = 40000
$ws.Range("M6").Value = -39887
$ws.Range("H7").Value = 433.33334
$ws.Range("I7").Value = 433.33334
$ws.Range("K7").Value = 433.33334
$ws.Range("M7").Value = -320.33334
$ws.Range("H31").Value = 2402.7144
$ws.Range("I31").Value = 1980
$ws.Range("K31").Value = 1980
$ws.Range("M31").Value = -1685
$ws.Range("H34").Value = 2402.7144
$ws.Range("I34").Value = 1980
$ws.Range("K34").Value = 1980
$ws.Range("M34").Value = -1778
$ws.Range("H99").Value = 1183
$ws.Range("I99").Value = 1174.5
$ws.Range("J99").Value = 1200
$ws.Range("K99").Value = 1174.5
$ws.Range("L99").Value = 1200
$ws.Range("M99").Value = 323.5
$ws.Range("N99").Value = -4196
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = ""
$ws.Range("H126").Value = 1183
$ws.Range("I126").Value = 1174.5
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 3523.5
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -1053.5
$ws.Range("N126").Value = -8540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 452.44446
$ws.Range("I4").Value = 383.75
$ws.Range("K4").Value = 1151.25
$ws.Range("M4").Value = -1039.25
$ws.Range("H68").Value = 936.9067
$ws.Range("I68").Value = 648.2857
$ws.Range("J68").Value = 1189.45
$ws.Range("K68").Value = 1944.8571
$ws.Range("L68").Value = 3568.35
$ws.Range("M68").Value = -1133.8571
$ws.Range("N68").Value = -5190.35
$ws.Range("H71").Value = 936.9067
$ws.Range("I71").Value = 648.2857
$ws.Range("J71").Value = 1189.45
$ws.Range("K71").Value = 5834.571300000001
$ws.Range("L71").Value = 10705.05
$ws.Range("M71").Value = -1778.571300000001
$ws.Range("N71").Value = -18817.05
$ws.Range("H131").Value = 1251726.1
$ws.Range("I131").Value = 1826.1666
$ws.Range("J131").Value = 1353069.4
$ws.Range("K131").Value = 5478.4998
$ws.Range("L131").Value = 4059208.2
$ws.Range("M131").Value = -438.4997999999996
$ws.Range("N131").Value = -4069288.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2827.4546
$ws.Range("J5").Value = 2827.4546
$ws.Range("L5").Value = 2827.4546
$ws.Range("N5").Value = -3051.4546
$ws.Range("H97").Value = 894.2
$ws.Range("I97").Value = 916.9167
$ws.Range("J97").Value = 349
$ws.Range("K97").Value = 916.9167
$ws.Range("L97").Value = 349
$ws.Range("M97").Value = -420.9167
$ws.Range("N97").Value = -1341
$ws.Range("H113").Value = 3396.8572
$ws.Range("I113").Value = 3155.6
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3155.6
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -985.5999999999999
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6540.8
$ws.Range("I40").Value = 6540.8
$ws.Range("K40").Value = 6540.8
$ws.Range("M40").Value = -6404.8
$ws.Range("H93").Value = 2857.7058
$ws.Range("I93").Value = 2173.5
$ws.Range("J93").Value = 4499.8
$ws.Range("K93").Value = 2173.5
$ws.Range("L93").Value = 4499.8
$ws.Range("M93").Value = -925.5
$ws.Range("N93").Value = -6995.8
$ws.Range("H122").Value = 21940.8
$ws.Range("I122").Value = 100004
$ws.Range("J122").Value = 2425
$ws.Range("K122").Value = 300012
$ws.Range("L122").Value = 7275
$ws.Range("M122").Value = -297562
$ws.Range("N122").Value = -12175

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 399.9524
$ws.Range("I113").Value = 300.1875
$ws.Range("J113").Value = 719.2
$ws.Range("K113").Value = 900.5625
$ws.Range("L113").Value = 2157.6
$ws.Range("M113").Value = 1269.4375
$ws.Range("N113").Value = -6497.6
$ws.Range("H126").Value = 945.875
$ws.Range("I126").Value = 919.8
$ws.Range("J126").Value = 1076.25
$ws.Range("K126").Value = 2759.4
$ws.Range("L126").Value = 3228.75
$ws.Range("M126").Value = -289.3999999999996
$ws.Range("N126").Value = -8168.75
$ws.Range("H132").Value = 4585.673
$ws.Range("I132").Value = 5560.189
$ws.Range("J132").Value = 2181.8667
$ws.Range("K132").Value = 16680.567
$ws.Range("L132").Value = 6545.6001
$ws.Range("M132").Value = -14150.567
$ws.Range("N132").Value = -11605.6001
